# Implemented Triangle Wave fit
# Adds 5 new ROI config rows (13-17) to Sheet1:
#   BMPDloop, NiLatticeDepthCalib, BMPDloopTof3000, PdBoBm, PDBO

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: BMPDloop
$ws.Cells.Item(13, 1).Value = "BMPDloop"
$ws.Cells.Item(13, 2).Value = 974
$ws.Cells.Item(13, 3).Value = 1554
$ws.Cells.Item(13, 4).Value = 1754
$ws.Cells.Item(13, 5).Value = 1830
$ws.Cells.Item(13, 6).Value = 2160
$ws.Cells.Item(13, 7).Value = 2560
$ws.Cells.Item(13, 8).Value = 349.1
$ws.Cells.Item(13, 9).Value = "[1064 1611 200 40]"
$ws.Cells.Item(13, 10).Value = "[2 1]"
$ws.Cells.Item(13, 11).Value = "[200 100]"

# Row 14: NiLatticeDepthCalib
$ws.Cells.Item(14, 1).Value = "NiLatticeDepthCalib"
$ws.Cells.Item(14, 2).Value = 1035
$ws.Cells.Item(14, 3).Value = 1491
$ws.Cells.Item(14, 4).Value = 1708
$ws.Cells.Item(14, 5).Value = 1838
$ws.Cells.Item(14, 6).Value = 2160
$ws.Cells.Item(14, 7).Value = 2560
$ws.Cells.Item(14, 8).Value = 349.1
$ws.Cells.Item(14, 9).Value = "[975 1588 50 50]"
$ws.Cells.Item(14, 10).Value = "[3 1]"
$ws.Cells.Item(14, 11).Value = "[130 100]"

# Row 15: BMPDloopTof3000
$ws.Cells.Item(15, 1).Value = "BMPDloopTof3000"
$ws.Cells.Item(15, 2).Value = 924
$ws.Cells.Item(15, 3).Value = 1604
$ws.Cells.Item(15, 4).Value = 1704
$ws.Cells.Item(15, 5).Value = 1880
$ws.Cells.Item(15, 6).Value = 2160
$ws.Cells.Item(15, 7).Value = 2560
$ws.Cells.Item(15, 8).Value = 349.1
$ws.Cells.Item(15, 9).Value = "[1015 1602 260 50]"
$ws.Cells.Item(15, 10).Value = "[2 1]"
$ws.Cells.Item(15, 11).Value = "[340 100]"

# Row 16: PdBoBm
$ws.Cells.Item(16, 1).Value = "PdBoBm"
$ws.Cells.Item(16, 2).Value = 994
$ws.Cells.Item(16, 3).Value = 1588
$ws.Cells.Item(16, 4).Value = 1746
$ws.Cells.Item(16, 5).Value = 1822
$ws.Cells.Item(16, 6).Value = 2160
$ws.Cells.Item(16, 7).Value = 2560
$ws.Cells.Item(16, 8).Value = 349.1
$ws.Cells.Item(16, 9).Value = "[]"
$ws.Cells.Item(16, 10).Value = "[1 1]"
$ws.Cells.Item(16, 11).Value = "[100 100]"

# Row 17: PDBO
$ws.Cells.Item(17, 1).Value = "PDBO"
$ws.Cells.Item(17, 2).Value = 1172
$ws.Cells.Item(17, 3).Value = 1332
$ws.Cells.Item(17, 4).Value = 1781
$ws.Cells.Item(17, 5).Value = 1819
$ws.Cells.Item(17, 6).Value = 2160
$ws.Cells.Item(17, 7).Value = 2560
$ws.Cells.Item(17, 8).Value = 349.1
$ws.Cells.Item(17, 9).Value = "[]"
$ws.Cells.Item(17, 10).Value = "[1 1]"
$ws.Cells.Item(17, 11).Value = "[100 100]"
